# Add two new rows (9 and 10) of results to the Summary sheet
# Bug fix: hydro inflow units corrected from GWh to MWh

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row9 = @(
    185851326104.2434,
    153358627530.0663,
    0,
    30461778814.17702,
    2030919760,
    183820406344.2434,
    2030919760,
    0,
    0,
    0,
    0,
    0,
    185851326104.2434,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    763.8449997901917,
    185851309924.4102,
    185851326104.2434,
    16179.83316040039,
    'optimal',
    'costs_emissionlimit',
    -1,
    -1,
    1,
    'GreenFieldHydro_Island_2040',
    'rawResults\20260203120741_GreenFieldHydro_Island_2040-1'
)

$row10 = @(
    168992919929.3019,
    143754608803.2541,
    0,
    27814901844.26778,
    2030919760,
    171569510647.5219,
    2030919760,
    2345674704.76664,
    -6953185182.986558,
    0,
    0,
    0,
    168992919929.3019,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    734.143000125885,
    168992917590.919,
    168992919929.3019,
    2338.382965087891,
    'optimal',
    'costs_emissionlimit',
    -1,
    -1,
    1,
    'GreenFieldHydro_NotIsland_2040',
    'rawResults\20260203173001_GreenFieldHydro_NotIsland_2040-1'
)

for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
}

for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, $i + 1).Value = $row10[$i]
}
